$d = $word.ActiveDocument

$pairs = @(
    @("396×7=", "750×6="),
    @("302×9=", "116×9="),
    @("700×8=", "674×4="),
    @("493×2=", "848×9="),
    @("636×7=", "703×6="),
    @("198×2=", "177×8="),
    @("137×2=", "343×6="),
    @("388×6=", "767×7="),
    @("672×2=", "669×3="),
    @("624×6=", "271×2="),
    @("350×8=", "455×3="),
    @("526×6=", "812×7="),
    @("761×7=", "933×9="),
    @("881×3=", "262×9="),
    @("173×5=", "300×9="),
    @("663×6=", "804×9="),
    @("338×6=", "720×3="),
    @("636×5=", "542×6="),
    @("937×3=", "149×6="),
    @("288×5=", "604×8="),
    @("636×8=", "403×7="),
    @("180×2=", "371×9="),
    @("268×8=", "984×7="),
    @("428×4=", "294×6="),
    @("222×5=", "949×5=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
